$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing cell H25: 4 -> 5 ---
$ws.Range("H25").Value = 5

# --- Row 35 (2020-04-17) ---
$ws.Range("A35").Value = 43938
$ws.Range("A35").NumberFormat = $ws.Range("A34").NumberFormat
$ws.Range("B35").Value = 12
$ws.Range("C35").Value = 450
$ws.Range("F35").Value = 246
$ws.Range("G35").Value = "Community(5)"
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0

# --- Row 36 (2020-04-18) ---
$ws.Range("A36").Value = 43939
$ws.Range("A36").NumberFormat = $ws.Range("A34").NumberFormat
$ws.Range("B36").Value = 16
$ws.Range("C36").Value = 1115
$ws.Range("F36").Value = 262
$ws.Range("G36").Value = "Community(16)"
$ws.Range("E36").Value = "Nairobi(9),Mombasa(5), Homabay(1)"
$ws.Range("H36").Value = 7
$ws.Range("I36").Value = 1
$ws.Range("L36").Value = "23-84"
$ws.Range("O36").Value = 4
$ws.Range("P36").Value = 12

# --- Row 37 (2020-04-19) ---
$ws.Range("A37").Value = 43940
$ws.Range("A37").NumberFormat = $ws.Range("A34").NumberFormat
$ws.Range("B37").Value = 8
$ws.Range("C37").Value = 1330
$ws.Range("F37").Value = 270
$ws.Range("G37").Value = "Community(8)"
$ws.Range("H37").Value = 7
$ws.Range("I37").Value = 2
$ws.Range("L37").Value = "17-65"

# --- Row 38 (2020-04-20) ---
$ws.Range("A38").Value = 43941
$ws.Range("A38").NumberFormat = $ws.Range("A34").NumberFormat
$ws.Range("B38").Value = 11
$ws.Range("E38").Value = "Mombasa(7),Nairobi(4)"
$ws.Range("F38").Value = 281
$ws.Range("G38").Value = "Community(11)"
$ws.Range("H38").Value = 2
$ws.Range("I38").Value = 0
$ws.Range("L38").Value = "11-80."
$ws.Range("L38").NumberFormat = $ws.Range("L31").NumberFormat
$ws.Range("O38").Value = 6
$ws.Range("P38").Value = 5

# --- Sheet view updates (scroll position + selection) ---
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("A38").Select() | Out-Null
